$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6 table: switch the applied table style.
#    {615EDEB4-4100-4870-A36F-BD64FAB464A8} -> {F5DB8DC3-18B3-4CD5-B0F2-87CE0E5E16EA}
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F5DB8DC3-18B3-4CD5-B0F2-87CE0E5E16EA}")
    }
}

# ---------------------------------------------------------------------------
# 2) Presentation theme (ppt/theme/theme1.xml, used by the slide master):
#    repaint the 12 theme colors from the "Integral" palette to the plain
#    "Office Theme" palette (dk1/lt1 are already shared, so only the other
#    ten entries actually change value).
#    Colors are exposed through Slide.ThemeColorScheme, items 1..12 map to
#    dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink (matching the
#    order of <a:clrScheme> children), and the RGB property takes/returns
#    a standard Win32 COLORREF (0x00BBGGRR encoded as R + G*256 + B*65536).
# ---------------------------------------------------------------------------
function ToCOLORREF($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ToCOLORREF($officeThemeColors[$i - 1])
}
